$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny precision revisions in existing cells
$ws.Range("E7").Value = 7.7251720428466797
$ws.Range("F9").Value = 24.59153938293457
$ws.Range("C11").Value = 14.144757270812988

# Add new row 23: lccNA_pcnt
$ws.Range("A23").Value = "lccNA_pcnt"
$ws.Range("B23").Value = 53.833560943603516
$ws.Range("C23").Value = 57.920494079589844
$ws.Range("D23").Value = 80.2430419921875
$ws.Range("E23").Value = 83.460487365722656
$ws.Range("F23").Value = 92.899467468261719
$ws.Range("G23").Value = 85.653961181640625
$ws.Range("H23").Value = 109.06099700927734

# Match formatting of the cells above (style index 1 = integer numFmt "0")
$ws.Range("B23:H23").NumberFormat = "0"
